# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (with per-fund holding detail, matching
# the layout of the existing "2021-Q3"/"2021-Q4" sheets) positioned right
# before the "总计" (totals) sheet, and prepends a corresponding summary row
# to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
#    NOTE: the worksheet object passed as the position anchor to .Add()
#    gets "taken over" by the newly created sheet in this COM shim, so we
#    must re-fetch "总计" by name afterwards rather than reuse that
#    variable.
# ---------------------------------------------------------------------
$q4Ws  = $wb.Worksheets.Item("2021-Q4")
$newWs = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newWs.Name = "2022-Q1"
$totalWs = $wb.Worksheets.Item("总计")

# Match the outline/grouping markers ("summaryBelow"/"summaryRight") used
# by every other sheet in this workbook.
$newWs.Outline.SummaryRow = 1
$newWs.Outline.SummaryColumn = 1

# Copy the header-row (B1:H1) + the index-column look (style only) from the
# "2021-Q4" sheet so the new sheet matches the existing visual style
# (bold, centered, thin-bordered header; bold centered index column).
$q4Ws.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$q4Ws.Range("A2").Copy()
$newWs.Range("A2:A6").PasteSpecial(-4122)

# Header labels.
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Fund holding rows. Columns B-G are text (even the numeric-looking ones,
# to match the source data / preserve leading zeros in fund codes);
# columns A (row index) and H (rank) are numbers.
$fundRows = @(
    @(0, "012368", "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金A", "5.76", "93.78", "6.76", "0.3894", 4),
    @(1, "233006", "大摩领先优势混合", "4.12", "94.42", "6.31", "0.2600", 5),
    @(2, "000309", "大摩品质生活精选股票", "4.36", "94.17", "5.93", "0.2585", 6),
    @(3, "010322", "摩根士丹利华鑫新兴产业股票", "2.41", "94.11", "6.49", "0.1564", 5),
    @(4, "012369", "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金C", "0.40", "93.78", "6.76", "0.0270", 4)
)

$r = 2
foreach ($row in $fundRows) {
    $newWs.Range("A$r").Value = $row[0]

    $textCols = @("B", "C", "D", "E", "F", "G")
    for ($i = 0; $i -lt $textCols.Length; $i++) {
        $addr = $textCols[$i] + $r
        $newWs.Range($addr).NumberFormat = "@"
        $newWs.Range($addr).Value = $row[$i + 1]
        $newWs.Range($addr).ClearFormats()
    }

    $newWs.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row to "总计", shifting older rows down.
# ---------------------------------------------------------------------
$totalWs.Rows.Item(2).Insert()

# Restore the index-column style (bold/centered) on the new A2 cell.
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("B2").ClearFormats()
$totalWs.Range("C2").Value = 5
$totalWs.Range("C2").ClearFormats()
$totalWs.Range("D2").Value = 1.09
$totalWs.Range("D2").ClearFormats()

# Renumber the index column for the rows that shifted down.
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 3. Restore the original active sheet/selection.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
$wb.Worksheets.Item("2021-Q3").Range("A1").Select()
